$d = $word.ActiveDocument

# The paragraph "...<m>vernis</m>.</ab>" has a trailing "." run, immediately
# after the "</m>" run and before the "</ab>" run, that must be deleted
# (while leaving the "</m>" and "</ab>" runs themselves untouched).
# "</m>." is unique across the document, so locate it directly.
$full = $d.Content.Text
$needle = "</m>."
$idx = $full.IndexOf($needle)

$dotStart = $idx + ($needle.Length - 1)
$dotEnd = $dotStart + 1

$dotRange = $d.Range($dotStart, $dotEnd)
$dotRange.Delete()
